# Outstandings.xlsx - "Add files via upload" edit
#
# Summary of the change being applied:
#   - A new purchase-invoice line (Sr. No 4) is appended to the
#     "Purchase 22-23" sheet: INV/23-24/1065 from Cassun Electricals,
#     dated 30-Sep-2023, outstanding amount 6080 (Bill amount = E36).
#   - The active sheet/selection bookkeeping moves: "Purchase 22-23"
#     becomes the selected tab (with the new row's area selected),
#     while "Sale 22-23" stops being the selected tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# --- Append the new invoice row (row 36) to "Purchase 22-23" ---------------
# Row 29 carries the same look (date/text/number cell styles + row height)
# the new row needs, so clone its formatting first, then overwrite values.
$ws1.Range("A29:F29").Copy()
$ws1.Range("A36:F36").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Rows("36:36").RowHeight = $ws1.Rows("29:29").RowHeight
$excel.CutCopyMode = $false

$ws1.Range("A36").Value = 4
$ws1.Range("B36").Value = 45199              # 30-Sep-2023
$ws1.Range("C36").Value = "INV/23-24/1065"
$ws1.Range("D36").Value = "Cassun Electricals"
$ws1.Range("E36").Value = 6080
$ws1.Range("F36").Formula = "=E36"

# --- View / selection state --------------------------------------------
# "Sale 22-23" selection moves to G33 and is no longer the active tab.
$ws2.Activate()
$ws2.Range("G33").Select()

# "Purchase 22-23" becomes the active tab, scrolled/selected near the new row.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A38").Select()
